# Update the "Förändrad" date column (C) for rows 2-15 from 2023-09-16 (45185)
# to 2023-10-05 (45204), matching the automatic update reflected in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
